# Add a new "2022-Q3" sheet, inserted right before the existing "2022-Q2" sheet.
$wb = $excel.ActiveWorkbook
$q2 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($q2)
$newSheet.Name = "2022-Q3"

# ---- Header row (copy the same headers used on the other quarterly sheets) ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$newSheet.Range("B1:H1").Copy()
$q2.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# ---- Data rows ----
# Columns B (fund code) and D..G (numeric-looking figures) are stored as plain
# text in this workbook (note the leading zeros on some fund codes), so force
# a text number format before assigning those values.
$textCols = @("B", "D", "E", "F", "G")
foreach ($col in $textCols) {
    $newSheet.Range($col + "2:" + $col + "7").NumberFormat = "@"
}

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "161030"
$newSheet.Range("C2").Value = "富国中证体育产业指数A"
$newSheet.Range("D2").Value = "1.59"
$newSheet.Range("E2").Value = "94.00"
$newSheet.Range("F2").Value = "4.52"
$newSheet.Range("G2").Value = "0.0719"
$newSheet.Range("H2").Value = 7

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "014179"
$newSheet.Range("C3").Value = "中银证券远见价值混合A"
$newSheet.Range("D3").Value = "1.56"
$newSheet.Range("E3").Value = "93.65"
$newSheet.Range("F3").Value = "3.96"
$newSheet.Range("G3").Value = "0.0618"
$newSheet.Range("H3").Value = 6

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "005571"
$newSheet.Range("C4").Value = "中银证券新能源灵活配置混合A"
$newSheet.Range("D4").Value = "0.53"
$newSheet.Range("E4").Value = "90.32"
$newSheet.Range("F4").Value = "5.04"
$newSheet.Range("G4").Value = "0.0267"
$newSheet.Range("H4").Value = 8

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "013278"
$newSheet.Range("C5").Value = "富国中证体育产业指数C"
$newSheet.Range("D5").Value = "0.42"
$newSheet.Range("E5").Value = "94.00"
$newSheet.Range("F5").Value = "4.52"
$newSheet.Range("G5").Value = "0.0190"
$newSheet.Range("H5").Value = 7

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "005572"
$newSheet.Range("C6").Value = "中银证券新能源灵活配置混合C"
$newSheet.Range("D6").Value = "0.25"
$newSheet.Range("E6").Value = "90.32"
$newSheet.Range("F6").Value = "5.04"
$newSheet.Range("G6").Value = "0.0126"
$newSheet.Range("H6").Value = 8

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "014180"
$newSheet.Range("C7").Value = "中银证券远见价值混合C"
$newSheet.Range("D7").Value = "0.16"
$newSheet.Range("E7").Value = "93.65"
$newSheet.Range("F7").Value = "3.96"
$newSheet.Range("G7").Value = "0.0063"
$newSheet.Range("H7").Value = 6

# Match the bold/centered/bordered header + index-column style used elsewhere.
$q2.Range("A2:A7").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)

# ---- Update the "总计" (summary) sheet: insert a new row for 2022-Q3 ----
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

# The values in column A are a 0-based running index; bump the pre-existing
# rows (now shifted down to rows 3:9) by 1 to keep that sequence consistent.
for ($r = 3; $r -le 9; $r++) {
    $cur = $summary.Range("A" + $r).Value2
    $summary.Range("A" + $r).Value = $cur + 1
}

$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.2

Write-Host "2022-Q3 sheet inserted and 总计 updated"
